# Reception and administrator section update:
#  - "Personas" sheet: new "fecha_salida" / "motivo_salida" columns, a
#    correction to an identificacion value, an updated num_acompanantes
#    count, and a new accompanying-person row (checkout of folio 1001 and
#    its household).
#  - "Encuestas" sheet: drop the two stray blank placeholder cells that
#    used to trail the first survey row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Personas
# ---------------------------------------------------------------------
$personas = $wb.Worksheets.Item("Personas")

# New header columns L (fecha_salida) and M (motivo_salida), matching the
# bold / bordered / centered look already used for the other headers.
$personas.Range("L1").Value = "fecha_salida"
$personas.Range("M1").Value = "motivo_salida"

$headerRange = $personas.Range("A1:M1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Row 2 (folio 1001, "wiliam ochoa", Titular): fix the identificacion typo,
# clear tutor_folio (titulares don't have a tutor), bump num_acompanantes
# now that a second acompanante checked in, and record the checkout info.
$personas.Range("C2").Value = 13123123128
$personas.Range("I2").ClearContents()
$personas.Range("K2").Value = 2
$personas.Range("L2").Value = "2025-12-17 12:29:16"
$personas.Range("M2").Value = "trabajo"

# Row 3 (1001-A, Gabriela de León): same checkout stamp.
$personas.Range("L3").Value = "2025-12-17 12:29:16"
$personas.Range("M3").Value = "trabajo"

# Row 4: new acompanante (1001-B, Luna Ochoa D L) tied to tutor 1001.
$personas.Range("A4").Value = "1001-B"
$personas.Range("B4").Value = "Luna Ochoa D L"
$personas.Range("C4").Value = "s123s123s12"
$personas.Range("D4").Value = 0
$personas.Range("E4").NumberFormat = "@"
$personas.Range("E4").Value = "2025-04-10"
$personas.Range("F4").Value = "Americana"
$personas.Range("G4").Value = "Femenino"
$personas.Range("H4").Value = "Acompañante"
$personas.Range("I4").Value = 1001
$personas.Range("J4").NumberFormat = "@"
$personas.Range("J4").Value = "2025-12-17"
$personas.Range("K4").Value = 0
$personas.Range("L4").Value = "2025-12-17 12:29:16"
$personas.Range("M4").Value = "trabajo"

[void]$personas.Range("A1").Select()

# ---------------------------------------------------------------------
# Encuestas
# ---------------------------------------------------------------------
$encuestas = $wb.Worksheets.Item("Encuestas")
$encuestas.Range("I2:J2").ClearContents()
